$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy style from an existing header cell (AB1) to the new header cells
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) # xlPasteFormats

# Fill data rows 2 through 41 with Wins=82, Losses=80, Ties=0
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 29).Value = 82  # AC
    $ws.Cells.Item($r, 30).Value = 80  # AD
    $ws.Cells.Item($r, 31).Value = 0   # AE
}
